# Generate Report for Handback
# Updates the "handoff"/"handback"/"HO xliff generate" timestamps for the
# row corresponding to file "1a127ea0-56a1-44cc-b973-905a48545852.md"
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# Row 3 = 1a127ea0-56a1-44cc-b973-905a48545852.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2017-02-09 08:32:07"

# --- zh-cn sheet -------------------------------------------------------
# Row 3 = 1a127ea0-56a1-44cc-b973-905a48545852.md
# Column H = "Correspond Handoff Datetime"
# Column L = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2017-02-09 08:31:47"
$wsZhCn.Range("L3").Value = "2017-02-09 08:32:40"

# --- de-de sheet -------------------------------------------------------
# Row 3 = 1a127ea0-56a1-44cc-b973-905a48545852.md
# Column H = "Correspond Handoff Datetime"
# Column L = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2017-02-09 08:32:07"
$wsDeDe.Range("L3").Value = "2017-02-09 08:33:06"
